# Applies the cryptos price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.918.48"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.304.52"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'305.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'97.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").Value = "'35.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'18.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("D13").Value = "'0.118"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").Value = "'6.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "2.661.33"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "2.296.87"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "42.829.67"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "'67.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "'236.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").Value = "'2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").Value = "'4.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'25.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").Value = "'167.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Value = "'9.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "'33.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'4.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  -2.52%  "
$ws.Range("D36").Value = "'17.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.32%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.0694"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("D43").Value = "2.008.04"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("D45").Value = "'18.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.08%  "
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").Value = "'2.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.42%  "
$ws.Range("D50").Value = "'54.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "2.528.10"
$ws.Range("E51").Value = "  -0.37%  "
